# Applies the edits described by the commit diff to Final_attendance_OC.docx
$d = $word.ActiveDocument

# 1) Department name in the "To," block.
$d.Content.Find.Execute("Department of Information Technology,", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Department of Computer Science and Engineering,", 2)

# 2) Letter date field.
$d.Content.Find.Execute(": 2020-06-25", $true, $false, $false, $false, $false,
                         $true, 1, $false, ": 2020-06-14", 2)

# 3) Subject line text ("hi" -> real subject). Whole-word match so we don't
#    clobber "Bharathi" elsewhere in the document.
$d.Content.Find.Execute("hi", $true, $true, $false, $false, $false,
                         $true, 1, $false, "Permission for Event Conduction", 2)

# 4) Team/club name. A single whole-word replace of "as" -> "Robotic Club"
#    fixes both the standalone "as" run and the "Team as" run (matched as
#    the whole word "as" inside that run), matching the diff's two hunks.
$d.Content.Find.Execute("as", $true, $true, $false, $false, $false,
                         $true, 1, $false, "Robotic Club", 2)

# 5) "work for <reason>" placeholder text.
$d.Content.Find.Execute("fds", $true, $true, $false, $false, $false,
                         $true, 1, $false, "fcdfhfb", 2)

# 6) Event start date (must run BEFORE the "to ..." replace below, since
#    that replace's new text would otherwise accidentally contain this
#    find text again).
$d.Content.Find.Execute(" 2020-06-17", $true, $false, $false, $false, $false,
                         $true, 1, $false, " 2020-06-15", 2)

# 7) Event end date.
$d.Content.Find.Execute(" to 2020-07-03", $true, $false, $false, $false, $false,
                         $true, 1, $false, " to 2020-06-17", 2)

# 8) Start hour.
$d.Content.Find.Execute(" from 3", $true, $false, $false, $false, $false,
                         $true, 1, $false, " from 8", 2)

# 9) Start minutes.
$d.Content.Find.Execute(":15", $true, $false, $false, $false, $false,
                         $true, 1, $false, ":30", 2)

# 10) End minutes.
$d.Content.Find.Execute(":50", $true, $false, $false, $false, $false,
                         $true, 1, $false, ":20", 2)

# 11) Event description placeholder paragraph ("fd" -> real sentence).
$d.Content.Find.Execute("fd", $true, $true, $false, $false, $false,
                         $true, 1, $false, "This is an workshop for 1st year Students.", 2)

# 12) Table row 1: student name and roll number placeholders.
$d.Content.Find.Execute("fdsf", $true, $true, $false, $false, $false,
                         $true, 1, $false, "Aaris", 2)
$d.Content.Find.Execute("ffd", $true, $true, $false, $false, $false,
                         $true, 1, $false, "18P61A05D7", 2)

# 13) Add a second student row to the attendance table: Yash / 18P61A05C2.
$t = $d.Tables.Item(1)
$newRow = $t.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "Yash"
$newRow.Cells.Item(2).Range.Text = "18P61A05C2"
